# Auto-generated script applying the scheduled market-price refresh diff
# to the Sagittarius_Profits workbook (columns H-N per Leve table).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 616.8333
$ws.Range("I15").Value = 616.8333
$ws.Range("K15").Value = 1850.4999
$ws.Range("M15").Value = -1681.4999
$ws.Range("H28").Value = 1247.5454
$ws.Range("J28").Value = 177
$ws.Range("L28").Value = 177
$ws.Range("N28").Value = -1147
$ws.Range("H33").Value = 164.125
$ws.Range("I33").Value = 164.125
$ws.Range("K33").Value = 164.125
$ws.Range("M33").Value = 64.875
$ws.Range("H43").Value = 0
$ws.Range("I43").Value = 0
$ws.Range("J43").Value = 0
$ws.Range("K43").Value = 0
$ws.Range("L43").Value = 0
$ws.Range("H52").Value = 1325
$ws.Range("I52").Value = 1325
$ws.Range("K52").Value = 3975
$ws.Range("M52").Value = -3815
$ws.Range("H80").Value = 2755.5
$ws.Range("J80").Value = 2737.4546
$ws.Range("L80").Value = 8212.363799999999
$ws.Range("N80").Value = -10208.3638
$ws.Range("H83").Value = 2755.5
$ws.Range("J83").Value = 2737.4546
$ws.Range("L83").Value = 24637.0914
$ws.Range("N83").Value = -34621.0914
$ws.Range("H86").Value = 4998
$ws.Range("I86").Value = 4998
$ws.Range("K86").Value = 4998
$ws.Range("M86").Value = -3875
$ws.Range("H89").Value = 4998
$ws.Range("I89").Value = 4998
$ws.Range("K89").Value = 24990
$ws.Range("M89").Value = -19374
$ws.Range("H116").Value = 4966.3335
$ws.Range("I116").Value = 4000
$ws.Range("J116").Value = 5449.5
$ws.Range("K116").Value = 4000
$ws.Range("L116").Value = 5449.5
$ws.Range("M116").Value = -558
$ws.Range("N116").Value = -12333.5
$ws.Range("H118").Value = 753
$ws.Range("J118").Value = 0
$ws.Range("L118").Value = 0
$ws.Range("H127").Value = 1967.75
$ws.Range("I127").Value = 2660
$ws.Range("K127").Value = 7980
$ws.Range("M127").Value = -3020
$ws.Range("H137").Value = 1569.85
$ws.Range("I137").Value = 1258.6471
$ws.Range("J137").Value = 3333.3333
$ws.Range("K137").Value = 3775.9413
$ws.Range("L137").Value = 9999.999899999999
$ws.Range("M137").Value = -1225.9413
$ws.Range("N137").Value = -15099.9999
$ws.Range("M43").ClearContents()
$ws.Range("N43").ClearContents()
$ws.Range("N118").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 664.6667
$ws.Range("I2").Value = 497
$ws.Range("K2").Value = 497
$ws.Range("M2").Value = -384
$ws.Range("H32").Value = 1751871.9
$ws.Range("J32").Value = 702044
$ws.Range("L32").Value = 702044
$ws.Range("N32").Value = -702618
$ws.Range("H105").Value = 100370
$ws.Range("J105").Value = 100370
$ws.Range("L105").Value = 100370
$ws.Range("N105").Value = -107358
$ws.Range("H116").Value = 664.6667
$ws.Range("I116").Value = 497
$ws.Range("K116").Value = 497
$ws.Range("M116").Value = 1797

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 664.6667
$ws.Range("I3").Value = 497
$ws.Range("K3").Value = 497
$ws.Range("M3").Value = -383
$ws.Range("H22").Value = 195.5
$ws.Range("I22").Value = 195.5
$ws.Range("K22").Value = 195.5
$ws.Range("M22").Value = -22.5
$ws.Range("H86").Value = 1466
$ws.Range("I86").Value = 1450
$ws.Range("J86").Value = 1498
$ws.Range("K86").Value = 1450
$ws.Range("L86").Value = 1498
$ws.Range("M86").Value = -327
$ws.Range("N86").Value = -3744
$ws.Range("H89").Value = 1466
$ws.Range("I89").Value = 1450
$ws.Range("J89").Value = 1498
$ws.Range("K89").Value = 7250
$ws.Range("L89").Value = 7490
$ws.Range("M89").Value = -1634
$ws.Range("N89").Value = -18722
$ws.Range("H132").Value = 30709
$ws.Range("I132").Value = 30709
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 30709
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -25649
$ws.Range("H134").Value = 2513.6667
$ws.Range("I134").Value = 2513.6667
$ws.Range("J134").Value = 0
$ws.Range("K134").Value = 7541.000100000001
$ws.Range("L134").Value = 0
$ws.Range("M134").Value = -5006.000100000001
$ws.Range("N132").ClearContents()
$ws.Range("N134").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H29").Value = 10000
$ws.Range("J29").Value = 10000
$ws.Range("L29").Value = 10000
$ws.Range("N29").Value = -10586
$ws.Range("H31").Value = 1311.6666
$ws.Range("I31").Value = 1272.2858
$ws.Range("K31").Value = 1272.2858
$ws.Range("M31").Value = -977.2858000000001
$ws.Range("H34").Value = 1311.6666
$ws.Range("I34").Value = 1272.2858
$ws.Range("K34").Value = 1272.2858
$ws.Range("M34").Value = -1070.2858
$ws.Range("H43").Value = 14828.5
$ws.Range("J43").Value = 14828.5
$ws.Range("L43").Value = 14828.5
$ws.Range("N43").Value = -15196.5
$ws.Range("H58").Value = 1673.8462
$ws.Range("I58").Value = 1628.1666
$ws.Range("J58").Value = 2222
$ws.Range("K58").Value = 1628.1666
$ws.Range("L58").Value = 2222
$ws.Range("M58").Value = -1425.1666
$ws.Range("N58").Value = -2628
$ws.Range("H94").Value = 163545.58
$ws.Range("J94").Value = 5273
$ws.Range("L94").Value = 5273
$ws.Range("N94").Value = -6175
$ws.Range("H101").Value = 14828.5
$ws.Range("J101").Value = 14828.5
$ws.Range("L101").Value = 14828.5
$ws.Range("N101").Value = -21318.5
$ws.Range("H107").Value = 1144.3334
$ws.Range("I107").Value = 1134.75
$ws.Range("K107").Value = 1134.75
$ws.Range("M107").Value = 785.25
$ws.Range("H108").Value = 0
$ws.Range("J108").Value = 0
$ws.Range("L108").Value = 0
$ws.Range("H132").Value = 3435.2727
$ws.Range("I132").Value = 3710.6667
$ws.Range("K132").Value = 11132.0001
$ws.Range("M132").Value = -8602.000100000001
$ws.Range("H136").Value = 1673.8462
$ws.Range("I136").Value = 1628.1666
$ws.Range("J136").Value = 2222
$ws.Range("K136").Value = 4884.4998
$ws.Range("L136").Value = 6666
$ws.Range("M136").Value = -2334.4998
$ws.Range("N136").Value = -11766
$ws.Range("N108").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 1581.8
$ws.Range("J68").Value = 1877.25
$ws.Range("L68").Value = 5631.75
$ws.Range("N68").Value = -7253.75
$ws.Range("H71").Value = 1581.8
$ws.Range("J71").Value = 1877.25
$ws.Range("L71").Value = 16895.25
$ws.Range("N71").Value = -25007.25
$ws.Range("H80").Value = 9333.333000000001
$ws.Range("I80").Value = 10000
$ws.Range("J80").Value = 9000
$ws.Range("K80").Value = 30000
$ws.Range("L80").Value = 27000
$ws.Range("M80").Value = -29064
$ws.Range("N80").Value = -28872
$ws.Range("H83").Value = 9333.333000000001
$ws.Range("I83").Value = 10000
$ws.Range("J83").Value = 9000
$ws.Range("K83").Value = 90000
$ws.Range("L83").Value = 81000
$ws.Range("M83").Value = -85320
$ws.Range("N83").Value = -90360

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H98").Value = 16100.556
$ws.Range("J98").Value = 16100.556
$ws.Range("L98").Value = 16100.556
$ws.Range("N98").Value = -22090.556
$ws.Range("H122").Value = 2224.8462
$ws.Range("I122").Value = 1866.25
$ws.Range("J122").Value = 2798.6
$ws.Range("K122").Value = 5598.75
$ws.Range("L122").Value = 8395.799999999999
$ws.Range("M122").Value = -3148.75
$ws.Range("N122").Value = -13295.8
$ws.Range("H140").Value = 142464.67
$ws.Range("J140").Value = 142464.67
$ws.Range("L140").Value = 142464.67
$ws.Range("N140").Value = -152824.67

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1249.25
$ws.Range("I22").Value = 1249.25
$ws.Range("J22").Value = 0
$ws.Range("K22").Value = 1249.25
$ws.Range("L22").Value = 0
$ws.Range("M22").Value = -954.25
$ws.Range("H27").Value = 1249.25
$ws.Range("I27").Value = 1249.25
$ws.Range("J27").Value = 0
$ws.Range("K27").Value = 1249.25
$ws.Range("L27").Value = 0
$ws.Range("M27").Value = -1142.25
$ws.Range("H31").Value = 437.66666
$ws.Range("I31").Value = 315
$ws.Range("J31").Value = 499
$ws.Range("K31").Value = 315
$ws.Range("L31").Value = 499
$ws.Range("M31").Value = -67
$ws.Range("N31").Value = -995
$ws.Range("H32").Value = 20000
$ws.Range("I32").Value = 20000
$ws.Range("K32").Value = 20000
$ws.Range("M32").Value = -19683
$ws.Range("H46").Value = 3664.889
$ws.Range("J46").Value = 4799.8
$ws.Range("L46").Value = 4799.8
$ws.Range("N46").Value = -5175.8
$ws.Range("H97").Value = 19895
$ws.Range("J97").Value = 19895
$ws.Range("L97").Value = 19895
$ws.Range("N97").Value = -21877
$ws.Range("H103").Value = 35000
$ws.Range("J103").Value = 35000
$ws.Range("L103").Value = 35000
$ws.Range("N103").Value = -37344
$ws.Range("N22").ClearContents()
$ws.Range("N27").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 2077.7144
$ws.Range("I122").Value = 2077.7144
$ws.Range("K122").Value = 6233.1432
$ws.Range("M122").Value = -3783.1432
$ws.Range("H132").Value = 1726.2632
$ws.Range("I132").Value = 1693.5294
$ws.Range("K132").Value = 5080.5882
$ws.Range("M132").Value = -2550.5882
$ws.Range("H141").Value = 181965.83
$ws.Range("J141").Value = 175559.2
$ws.Range("L141").Value = 175559.2
